# Auto-generated Excel COM-interop script applying the Seraph_Profits.xlsx diff.
# Updates specific H/I/J/K/L/M/N numeric cells across the ALC, ARM, BSM, CRP, CUL,
# GSM, LTW and WVR worksheets (current price / profit recalculation values).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 76923440
$ws.Range("I28").Value = 76923440
$ws.Range("K28").Value = 76923440
$ws.Range("M28").Value = -76922955
$ws.Range("H86").Value = 1904
$ws.Range("J86").Value = 1904
$ws.Range("L86").Value = 1904
$ws.Range("N86").Value = -4150
$ws.Range("H89").Value = 1904
$ws.Range("J89").Value = 1904
$ws.Range("L89").Value = 9520
$ws.Range("N89").Value = -20752

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 40.25
$ws.Range("I5").Value = 33.142857
$ws.Range("K5").Value = 33.142857
$ws.Range("M5").Value = 78.85714300000001
$ws.Range("H32").Value = 5935.517
$ws.Range("I32").Value = 4135.6523
$ws.Range("K32").Value = 4135.6523
$ws.Range("M32").Value = -3848.6523
$ws.Range("H74").Value = 16006
$ws.Range("I74").Value = 7012
$ws.Range("K74").Value = 7012
$ws.Range("M74").Value = -6138
$ws.Range("H77").Value = 16006
$ws.Range("I77").Value = 7012
$ws.Range("K77").Value = 35060
$ws.Range("M77").Value = -30692

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 40.25
$ws.Range("I4").Value = 33.142857
$ws.Range("K4").Value = 33.142857
$ws.Range("M4").Value = 81.85714300000001
$ws.Range("H22").Value = 3661.3333
$ws.Range("I22").Value = 492
$ws.Range("K22").Value = 492
$ws.Range("M22").Value = -319
$ws.Range("H105").Value = 4632304
$ws.Range("I105").Value = 10418456
$ws.Range("J105").Value = 3383.1
$ws.Range("K105").Value = 10418456
$ws.Range("L105").Value = 3383.1
$ws.Range("M105").Value = -10416709
$ws.Range("N105").Value = -6877.1

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5265.3
$ws.Range("I31").Value = 4564.3076
$ws.Range("K31").Value = 4564.3076
$ws.Range("M31").Value = -4269.3076
$ws.Range("H34").Value = 5265.3
$ws.Range("I34").Value = 4564.3076
$ws.Range("K34").Value = 4564.3076
$ws.Range("M34").Value = -4362.3076
$ws.Range("H58").Value = 4724.4
$ws.Range("I58").Value = 3316.5
$ws.Range("K58").Value = 3316.5
$ws.Range("M58").Value = -3113.5
$ws.Range("H99").Value = 14885.333
$ws.Range("I99").Value = 9961.9
$ws.Range("J99").Value = 19361.182
$ws.Range("K99").Value = 9961.9
$ws.Range("L99").Value = 19361.182
$ws.Range("M99").Value = -8463.9
$ws.Range("N99").Value = -22357.182
$ws.Range("H126").Value = 14885.333
$ws.Range("I126").Value = 9961.9
$ws.Range("J126").Value = 19361.182
$ws.Range("K126").Value = 29885.7
$ws.Range("L126").Value = 58083.546
$ws.Range("M126").Value = -27415.7
$ws.Range("N126").Value = -63023.546
$ws.Range("H132").Value = 2321.375
$ws.Range("I132").Value = 2321.375
$ws.Range("K132").Value = 6964.125
$ws.Range("M132").Value = -4434.125
$ws.Range("H134").Value = 3542.875
$ws.Range("I134").Value = 3117.3333
$ws.Range("K134").Value = 9351.999899999999
$ws.Range("M134").Value = -6816.999899999999
$ws.Range("H136").Value = 4724.4
$ws.Range("I136").Value = 3316.5
$ws.Range("K136").Value = 9949.5
$ws.Range("M136").Value = -7399.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1874.75
$ws.Range("I11").Value = 1874.75
$ws.Range("K11").Value = 5624.25
$ws.Range("M11").Value = -5484.25
$ws.Range("H26").Value = 332.8
$ws.Range("J26").Value = 507.6
$ws.Range("L26").Value = 1522.8
$ws.Range("N26").Value = -2098.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H57").Value = 11698.833
$ws.Range("I57").Value = 7000
$ws.Range("J57").Value = 14048.25
$ws.Range("K57").Value = 7000
$ws.Range("L57").Value = 14048.25
$ws.Range("M57").Value = -6180
$ws.Range("N57").Value = -15688.25
$ws.Range("H80").Value = 3725
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 5000
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 3725
$ws.Range("J83").Value = 5000
$ws.Range("L83").Value = 25000
$ws.Range("N83").Value = -34984

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1047
$ws.Range("I22").Value = 1146
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 1146
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -851
$ws.Range("N22").Value = -1340
$ws.Range("H27").Value = 1047
$ws.Range("I27").Value = 1146
$ws.Range("J27").Value = 750
$ws.Range("K27").Value = 1146
$ws.Range("L27").Value = 750
$ws.Range("M27").Value = -1039
$ws.Range("N27").Value = -964
$ws.Range("H40").Value = 1500
$ws.Range("I40").Value = 1500
$ws.Range("K40").Value = 1500
$ws.Range("M40").Value = -1364
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H61").Value = 13890818
$ws.Range("I61").Value = 18519518
$ws.Range("K61").Value = 18519518
$ws.Range("M61").Value = -18519316
$ws.Range("H69").Value = 55000
$ws.Range("I69").Value = 55000
$ws.Range("K69").Value = 55000
$ws.Range("M69").Value = -54189
$ws.Range("H72").Value = 55000
$ws.Range("I72").Value = 55000
$ws.Range("K72").Value = 165000
$ws.Range("M72").Value = -160944
$ws.Range("H93").Value = 1121.1
$ws.Range("I93").Value = 899.1539
$ws.Range("K93").Value = 899.1539
$ws.Range("M93").Value = 348.8461
$ws.Range("H113").Value = 13890818
$ws.Range("I113").Value = 18519518
$ws.Range("K113").Value = 18519518
$ws.Range("M113").Value = -18517348
$ws.Range("H132").Value = 1584.8572
$ws.Range("I132").Value = 1639
$ws.Range("J132").Value = 1449.5
$ws.Range("K132").Value = 4917
$ws.Range("L132").Value = 4348.5
$ws.Range("M132").Value = -2387
$ws.Range("N132").Value = -9408.5
$ws.Range("H136").Value = 5606.4
$ws.Range("I136").Value = 5372
$ws.Range("J136").Value = 8888
$ws.Range("K136").Value = 16116
$ws.Range("L136").Value = 26664
$ws.Range("M136").Value = -13566
$ws.Range("N136").Value = -31764

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H126").Value = 2179.6365
$ws.Range("I126").Value = 1872.75
$ws.Range("K126").Value = 5618.25
$ws.Range("M126").Value = -3148.25
$ws.Range("H132").Value = 3399.5
$ws.Range("I132").Value = 2974.3333
$ws.Range("K132").Value = 8922.999899999999
$ws.Range("M132").Value = -6392.999899999999

